$d = $word.ActiveDocument

# -----------------------------------------------------------------
# Edit 1: "Objetivos" paragraph - split the run into two <w:t> runs
# joined by a manual line break (<w:br/>) right where the original
# text ran "poliméricos." straight into "Fornecer" with no space.
# -----------------------------------------------------------------
$find1 = "materiais poliméricos.Fornecer conhecimentos técnicos"
$repl1 = "materiais poliméricos.^lFornecer conhecimentos técnicos"
$r1 = $d.Content
$found1 = $r1.Find.Execute($find1, $true, $false, $false, $false, $false, $true, 1, $false, $repl1, 2)
if (-not $found1) {
    throw "Edit 1 (Objetivos): target text not found"
}

# -----------------------------------------------------------------
# Edit 2: "Norma de recuperação" paragraph - insert two manual line
# breaks between "...pela fórmula: " and "MR = (NF + PR)/2".
# -----------------------------------------------------------------
$find2 = "calculada pela fórmula: MR = (NF + PR)/2"
$repl2 = "calculada pela fórmula: ^l^lMR = (NF + PR)/2"
$r2 = $d.Content
$found2 = $r2.Find.Execute($find2, $true, $false, $false, $false, $false, $true, 1, $false, $repl2, 2)
if (-not $found2) {
    throw "Edit 2 (Norma de recuperação) : target text not found"
}

# -----------------------------------------------------------------
# Edit 3: "Bibliografia" paragraph - insert a manual line break
# between every reference that used to run together with no
# separator.
# -----------------------------------------------------------------
$find3 = "BRETAS, R. E. S.; D´ÁVILA, M. A. Reologia de Polímeros Fundidos, São Carlos, Eduscar, 2005.MANRICH, S. Processamento de termoplásticos – Rosca única, extrusão & matrizes, injeção & moldes,. McCRUM, N. G., BUCKLEY, C. P., BUCKNALl, C. B. Principles of Polymer Engineering, New York, Oxford University Press, 1997.Blass A., Processamento de Polímeros, editora da UFSC.CHAWLA, K. K. Composite Materials Science and Engineering, Spring-Verlag ed., Berlin, 1987.BRETT, A.M.O., BRETT, C.M. Electroquímica: Princípios, métodos e aplicações. Livraria Medina, Coimbra, 1996.FONTANA, M. G. Corrosion Engineering. 3ª Edição. McGraw-Hill, 1987GENTIL, V. Corrosão. 5ª Edição, Rio de Janeiro, Ed. LTC, 2007 RAMANHATAN, L. Corrosão e seu Controle. São Paulo. Ed. Hemus, 1990SHREIR, L.L., JARMAN, R.A., BURSTEIN, G.T. Corrosion. 3ª Edição. Oxford, Butterworth Heinemann, volume 2, 2000WOLYNEC, S. Técnicas Eletroquímicas em Corrosão, EDUSP, São Paulo, 2003"

$repl3 = "BRETAS, R. E. S.; D´ÁVILA, M. A. Reologia de Polímeros Fundidos, São Carlos, Eduscar, 2005.^lMANRICH, S. Processamento de termoplásticos – Rosca única, extrusão & matrizes, injeção & moldes,. ^lMcCRUM, N. G., BUCKLEY, C. P., BUCKNALl, C. B. Principles of Polymer Engineering, New York, Oxford University Press, 1997.^lBlass A., Processamento de Polímeros, editora da UFSC.^lCHAWLA, K. K. Composite Materials Science and Engineering, Spring-Verlag ed., Berlin, 1987.^lBRETT, A.M.O., BRETT, C.M. Electroquímica: Princípios, métodos e aplicações. Livraria Medina, Coimbra, 1996.^lFONTANA, M. G. Corrosion Engineering. 3ª Edição. McGraw-Hill, 1987^lGENTIL, V. Corrosão. 5ª Edição, Rio de Janeiro, Ed. LTC, 2007 ^lRAMANHATAN, L. Corrosão e seu Controle. São Paulo. Ed. Hemus, 1990^lSHREIR, L.L., JARMAN, R.A., BURSTEIN, G.T. Corrosion. 3ª Edição. Oxford, Butterworth Heinemann, volume 2, 2000^lWOLYNEC, S. Técnicas Eletroquímicas em Corrosão, EDUSP, São Paulo, 2003"

$r3 = $d.Content
$found3 = $r3.Find.Execute($find3, $true, $false, $false, $false, $false, $true, 1, $false, $repl3, 2)
if (-not $found3) {
    throw "Edit 3 (Bibliografia): target text not found"
}

Write-Output "Edit1=$found1 Edit2=$found2 Edit3=$found3"
